# Insert a new data row at row 237 (pushing the existing rows 237-324 down
# to 238-325) and populate it with a new daily price observation, matching
# the commit "Fruta / hortaliza, semanal" (weekly fruit/vegetable price
# update for Puerro / Vega Modelo de Temuco).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(237).Insert()

$ws.Range("A237").Value = 10
$ws.Range("B237").Value = "Vega Modelo de Temuco"
$ws.Range("C237").Value = "La Araucanía"
$ws.Range("D237").Value = 45146
$ws.Range("E237").Value = 9
$ws.Range("F237").Value = 100112005
$ws.Range("G237").Value = "Puerro"
$ws.Range("H237").Value = "Azul de Maquehue"
$ws.Range("I237").Value = "Primera"
$ws.Range("J237").Value = 50
$ws.Range("K237").Value = 8000
$ws.Range("L237").Value = 8000
$ws.Range("M237").Value = 8000
$ws.Range("N237").Value = "$/docena de paquetes"
$ws.Range("O237").Value = "Provincia de Cautín"
$ws.Range("P237").Value = 667
$ws.Range("Q237").Value = 12
$ws.Range("R237").Value = "Hortaliza"
